$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 17882710
$ws.Range("J40").Value = 41706190
$ws.Range("L40").Value = 41706190
$ws.Range("N40").Value = -41706540
$ws.Range("H135").Value = 5279.8
$ws.Range("J135").Value = 7739
$ws.Range("L135").Value = 69651
$ws.Range("N135").Value = -74721
$ws.Range("H140").Value = 70344.14
$ws.Range("J140").Value = 70283.336
$ws.Range("L140").Value = 70283.336
$ws.Range("N140").Value = -80643.336
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 16299.357
$ws.Range("I61").Value = 19681.182
$ws.Range("K61").Value = 19681.182
$ws.Range("M61").Value = -19469.182
$ws.Range("H63").Value = 1351
$ws.Range("I63").Value = 1468
$ws.Range("J63").Value = 1000
$ws.Range("K63").Value = 1468
$ws.Range("L63").Value = 1000
$ws.Range("M63").Value = -782
$ws.Range("N63").Value = -2372
$ws.Range("H66").Value = 1351
$ws.Range("I66").Value = 1468
$ws.Range("J66").Value = 1000
$ws.Range("K66").Value = 7340
$ws.Range("L66").Value = 5000
$ws.Range("M66").Value = -3908
$ws.Range("N66").Value = -11864
$ws.Range("H74").Value = 13890350
$ws.Range("I74").Value = 83334170
$ws.Range("J74").Value = 1586.8667
$ws.Range("K74").Value = 83334170
$ws.Range("L74").Value = 1586.8667
$ws.Range("M74").Value = -83333296
$ws.Range("N74").Value = -3334.8667
$ws.Range("H77").Value = 13890350
$ws.Range("I77").Value = 83334170
$ws.Range("J77").Value = 1586.8667
$ws.Range("K77").Value = 416670850
$ws.Range("L77").Value = 7934.333500000001
$ws.Range("M77").Value = -416666482
$ws.Range("N77").Value = -16670.3335
$ws.Range("H97").Value = 658.5714
$ws.Range("I97").Value = 658.5714
$ws.Range("K97").Value = 658.5714
$ws.Range("M97").Value = -162.5714
$ws.Range("H132").Value = 12468.519
$ws.Range("I132").Value = 14290.721
$ws.Range("K132").Value = 42872.163
$ws.Range("M132").Value = -40342.163
$ws.Range("H136").Value = 16299.357
$ws.Range("I136").Value = 19681.182
$ws.Range("K136").Value = 59043.546
$ws.Range("M136").Value = -56493.546
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 25001440
$ws.Range("I105").Value = 34484052
$ws.Range("K105").Value = 34484052
$ws.Range("M105").Value = -34482305
$ws.Range("H134").Value = 2740.5557
$ws.Range("I134").Value = 1150.8
$ws.Range("K134").Value = 3452.4
$ws.Range("M134").Value = -917.3999999999996
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1292.6666
$ws.Range("J22").Value = 1781.75
$ws.Range("L22").Value = 1781.75
$ws.Range("N22").Value = -2481.75
$ws.Range("H31").Value = 4539.4365
$ws.Range("I31").Value = 845.9
$ws.Range("J31").Value = 6650.029
$ws.Range("K31").Value = 845.9
$ws.Range("L31").Value = 6650.029
$ws.Range("M31").Value = -550.9
$ws.Range("N31").Value = -7240.029
$ws.Range("H34").Value = 4539.4365
$ws.Range("I34").Value = 845.9
$ws.Range("J34").Value = 6650.029
$ws.Range("K34").Value = 845.9
$ws.Range("L34").Value = 6650.029
$ws.Range("M34").Value = -643.9
$ws.Range("N34").Value = -7054.029
$ws.Range("H58").Value = 375336.94
$ws.Range("I58").Value = 1252513.1
$ws.Range("K58").Value = 1252513.1
$ws.Range("M58").Value = -1252310.1
$ws.Range("H99").Value = 11227.24
$ws.Range("I99").Value = 13425.909
$ws.Range("J99").Value = 9499.714
$ws.Range("K99").Value = 13425.909
$ws.Range("L99").Value = 9499.714
$ws.Range("M99").Value = -11927.909
$ws.Range("N99").Value = -12495.714
$ws.Range("H126").Value = 11227.24
$ws.Range("I126").Value = 13425.909
$ws.Range("J126").Value = 9499.714
$ws.Range("K126").Value = 40277.727
$ws.Range("L126").Value = 28499.142
$ws.Range("M126").Value = -37807.727
$ws.Range("N126").Value = -33439.142
$ws.Range("H132").Value = 55591852
$ws.Range("J132").Value = 2999.5
$ws.Range("L132").Value = 8998.5
$ws.Range("N132").Value = -14058.5
$ws.Range("H134").Value = 1911.1351
$ws.Range("I134").Value = 1467.44
$ws.Range("J134").Value = 2835.5
$ws.Range("K134").Value = 4402.32
$ws.Range("L134").Value = 8506.5
$ws.Range("M134").Value = -1867.32
$ws.Range("N134").Value = -13576.5
$ws.Range("H136").Value = 375336.94
$ws.Range("I136").Value = 1252513.1
$ws.Range("K136").Value = 3757539.3
$ws.Range("M136").Value = -3754989.3
$ws.Range("H141").Value = 84226.44500000001
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 712.6923
$ws.Range("I5").Value = 717.8
$ws.Range("J5").Value = 695.6667
$ws.Range("K5").Value = 2153.4
$ws.Range("L5").Value = 2087.0001
$ws.Range("M5").Value = -2041.4
$ws.Range("N5").Value = -2311.0001
$ws.Range("H22").Value = 783.6667
$ws.Range("I22").Value = 249
$ws.Range("K22").Value = 747
$ws.Range("M22").Value = -578
$ws.Range("H27").Value = 783.6667
$ws.Range("I27").Value = 249
$ws.Range("K27").Value = 747
$ws.Range("M27").Value = -645
$ws.Range("H86").Value = 307.5
$ws.Range("I86").Value = 315
$ws.Range("J86").Value = 300
$ws.Range("K86").Value = 945
$ws.Range("L86").Value = 900
$ws.Range("M86").Value = 241
$ws.Range("N86").Value = -3272
$ws.Range("H89").Value = 307.5
$ws.Range("I89").Value = 315
$ws.Range("J89").Value = 300
$ws.Range("K89").Value = 2835
$ws.Range("L89").Value = 2700
$ws.Range("M89").Value = 3093
$ws.Range("N89").Value = -14556
$ws.Range("H135").Value = 712.6923
$ws.Range("I135").Value = 717.8
$ws.Range("J135").Value = 695.6667
$ws.Range("K135").Value = 6460.2
$ws.Range("L135").Value = 6261.0003
$ws.Range("M135").Value = -3925.2
$ws.Range("N135").Value = -11331.0003
$ws.Range("H137").Value = 57706016
$ws.Range("I137").Value = 93751930
$ws.Range("J137").Value = 32560
$ws.Range("K137").Value = 281255790
$ws.Range("L137").Value = 97680
$ws.Range("M137").Value = -281250690
$ws.Range("N137").Value = -107880
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2386145.5
$ws.Range("I70").Value = 2981304.2
$ws.Range("K70").Value = 2981304.2
$ws.Range("M70").Value = -2981034.2
$ws.Range("H73").Value = 2386145.5
$ws.Range("I73").Value = 2981304.2
$ws.Range("K73").Value = 2981304.2
$ws.Range("M73").Value = -2980368.2
$ws.Range("H122").Value = 443367.1
$ws.Range("I122").Value = 649687.4399999999
$ws.Range("J122").Value = 4936.375
$ws.Range("K122").Value = 1949062.32
$ws.Range("L122").Value = 14809.125
$ws.Range("M122").Value = -1946612.32
$ws.Range("N122").Value = -19709.125
$ws.Range("H132").Value = 47075.66
$ws.Range("I132").Value = 64364.695
$ws.Range("K132").Value = 193094.085
$ws.Range("M132").Value = -190564.085
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4344.6787
$ws.Range("I7").Value = 2868.9
$ws.Range("J7").Value = 5164.5557
$ws.Range("K7").Value = 2868.9
$ws.Range("L7").Value = 5164.5557
$ws.Range("M7").Value = -2756.9
$ws.Range("N7").Value = -5388.5557
$ws.Range("H40").Value = 30305450
$ws.Range("I40").Value = 937.1667
$ws.Range("K40").Value = 937.1667
$ws.Range("M40").Value = -801.1667
$ws.Range("H122").Value = 36871160
$ws.Range("I122").Value = 52635280
$ws.Range("K122").Value = 157905840
$ws.Range("M122").Value = -157903390
$ws.Range("H126").Value = 4344.6787
$ws.Range("I126").Value = 2868.9
$ws.Range("J126").Value = 5164.5557
$ws.Range("K126").Value = 8606.700000000001
$ws.Range("L126").Value = 15493.6671
$ws.Range("M126").Value = -6136.700000000001
$ws.Range("N126").Value = -20433.6671
$ws.Range("H132").Value = 3735.255
$ws.Range("I132").Value = 2938.919
$ws.Range("J132").Value = 5839.857
$ws.Range("K132").Value = 8816.757
$ws.Range("L132").Value = 17519.571
$ws.Range("M132").Value = -6286.757
$ws.Range("N132").Value = -22579.571
$ws.Range("H136").Value = 3657.9473
$ws.Range("I136").Value = 2850
$ws.Range("K136").Value = 8550
$ws.Range("M136").Value = -6000
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4782.95
$ws.Range("I122").Value = 4315.185
$ws.Range("K122").Value = 12945.555
$ws.Range("M122").Value = -10495.555
$ws.Range("H126").Value = 3472.0557
$ws.Range("I126").Value = 2935.7
$ws.Range("J126").Value = 4142.5
$ws.Range("K126").Value = 8807.099999999999
$ws.Range("L126").Value = 12427.5
$ws.Range("M126").Value = -6337.099999999999
$ws.Range("N126").Value = -17367.5
$ws.Range("H136").Value = 8512.492
$ws.Range("I136").Value = 2365.5217
$ws.Range("K136").Value = 7096.5651
$ws.Range("M136").Value = -4546.5651
